$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: human-readable column headers (replacing slug-style names)
$ws.Cells.Item(1,1).Value = "Edad (grandes grupos)"
$ws.Cells.Item(1,2).Value = "Personas"
$ws.Cells.Item(1,3).Value = "Residencia provincia código"
$ws.Cells.Item(1,4).Value = "Nacimiento provincia código"
$ws.Cells.Item(1,5).Value = "Residencia CCAA nombre"
$ws.Cells.Item(1,6).Value = "Provincia de nacimiento"
$ws.Cells.Item(1,7).Value = "Residencia provincia nombre"
$ws.Cells.Item(1,8).Value = "Año"
$ws.Cells.Item(1,9).Value = "Sexo"

# Row 2: measure/dimension identifiers
$ws.Cells.Item(2,1).Value = "iaest-measure:edad-grandes-grupos"
$ws.Cells.Item(2,2).Value = "iaest-measure:personas"
$ws.Cells.Item(2,3).Value = "null"
$ws.Cells.Item(2,4).Value = "null"
$ws.Cells.Item(2,5).Value = "iaest-measure:residencia-ccaa-nombre"
$ws.Cells.Item(2,6).Value = "iaest-measure:provincia-de-nacimiento"
$ws.Cells.Item(2,7).Value = "iaest-measure:residencia-provincia-nombre"
$ws.Cells.Item(2,8).Value = "sdmx-dimension:refPeriod"
$ws.Cells.Item(2,9).Value = "iaest-measure:sexo"

# Row 3: dim / medida classification
$ws.Cells.Item(3,1).Value = "medida"
$ws.Cells.Item(3,2).Value = "medida"
$ws.Cells.Item(3,3).Value = "null"
$ws.Cells.Item(3,4).Value = "null"
$ws.Cells.Item(3,5).Value = "medida"
$ws.Cells.Item(3,6).Value = "medida"
$ws.Cells.Item(3,7).Value = "medida"
$ws.Cells.Item(3,8).Value = "dim"
$ws.Cells.Item(3,9).Value = "medida"

# Row 4: datatypes (xsd:string replaces skos:Concept, xsd:date added for refPeriod)
$ws.Cells.Item(4,1).Value = "xsd:string"
$ws.Cells.Item(4,2).Value = "xsd:int"
$ws.Cells.Item(4,3).Value = "null"
$ws.Cells.Item(4,4).Value = "null"
$ws.Cells.Item(4,5).Value = "xsd:string"
$ws.Cells.Item(4,6).Value = "xsd:string"
$ws.Cells.Item(4,7).Value = "xsd:string"
$ws.Cells.Item(4,8).Value = "xsd:date"
$ws.Cells.Item(4,9).Value = "xsd:string"

# Row 5: mapping file moves from column A (ano) to column H (Año)
$ws.Cells.Item(5,1).Clear()
$ws.Cells.Item(5,8).Value = "mapping-ano.xlsx"
$ws.Cells.Item(4,8).Copy()
$ws.Cells.Item(5,8).PasteSpecial(-4122)
